$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 columns B:E
$ws.Range("B2").Value = 49.443367045917803
$ws.Range("C2").Value = 23.218118119056388
$ws.Range("D2").Value = 58.410035592184492
$ws.Range("E2").Value = 20.821599490754718

# Update row 3 columns B:E
$ws.Range("B3").Value = 46.028977094461943
$ws.Range("C3").Value = 21.561630270302388
$ws.Range("D3").Value = 84.036402381209641
$ws.Range("E3").Value = 32.706277011313034

# Update selection to reflect new selected range B1:E3
$ws.Range("B1:E3").Select()
